$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, 0.8514969533240468, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, 1.764146564639145, -15.73224693708785, 3.685821837315178, -15.73224693708785, 1.81564921395372, -15.73224693708785, 2.512827552082161, -15.73224693708785),
    @(-15.73224693708785, 4.321902916743732, -15.73224693708785, -15.73224693708785, -15.73224693708785, 3.25429788017431, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(2.894748942750159, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, 1.840431087544417, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(3.651054460154551, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, 0.7951842509482159, -15.73224693708785, 2.253578326660847),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, 1.947789735200385, -15.73224693708785, 1.955736940690373, -15.73224693708785, -15.73224693708785, -15.73224693708785, 1.508739995165169),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, 1.863192676610178, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, 1.139718443688061, 2.378506526498993),
    @(-15.73224693708785, -15.73224693708785, 1.405274597823634, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, 1.769131559563189),
    @(-15.73224693708785, -15.73224693708785, 0.7586822631344369, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, 2.051294892257054, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, 1.033932539431976, -15.73224693708785, -15.73224693708785, -15.73224693708785, 2.06883616902181, 0.6641738833969648, 1.978915876168003, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, -15.73224693708785, 2.007540263930589, 0.4046841148345513, 2.000801495425572, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, 1.937862942728536, -15.73224693708785, -15.73224693708785, -15.73224693708785, 1.825403730745636, 1.413058350177869, -15.73224693708785, -15.73224693708785),
    @(-15.73224693708785, -15.73224693708785, 2.674327631860859, -15.73224693708785, 2.834048709668738, -15.73224693708785, 1.209454081928172, 3.444746710189647, -15.73224693708785, 1.917240818862911),
    @(-15.73224693708785, -15.73224693708785, -15.73224693708785, 3.157461994049395, -15.73224693708785, 2.717760084750838, 1.266971362008894, -15.73224693708785, -15.73224693708785, -15.73224693708785)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
